$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date cell's formatting (style s="1", numFmtId 14) onto the new
# date cells so no new cellXfs entry gets created.
$ws.Range("A2").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122)

# Row 5 - Stage 3 results (2026-01-23)
$ws.Range("A5").Value = 46045
$ws.Range("B5").Value = "Santos Tour Down Under"
$ws.Range("C5").Value = "Stage 3"
$ws.Range("D5").Value = "Sam Welsford"
$ws.Range("E5").Value = "Tobias Lund Andresen"
$ws.Range("F5").Value = "Lewis Bower"
$ws.Range("G5").Value = "Jake Stewart"
$ws.Range("H5").Value = "Aaron Gate"
$ws.Range("I5").Value = "Žak Eržen"
$ws.Range("J5").Value = "Finn Fisher-Black"
$ws.Range("K5").Value = "Matthew Fox"
$ws.Range("L5").Value = "Anthon Charmig"
$ws.Range("M5").Value = "Jensen Plowright"

# Row 6 - Stage 4 results (2026-01-24)
$ws.Range("A6").Value = 46046
$ws.Range("B6").Value = "Santos Tour Down Under"
$ws.Range("C6").Value = "Stage 4"
$ws.Range("D6").Value = "Ethan Vernon"
$ws.Range("E6").Value = "Tobias Lund Andresen"
$ws.Range("F6").Value = "Laurence Pithie"
$ws.Range("G6").Value = "Brady Gilmore"
$ws.Range("H6").Value = "Aaron Gate"
$ws.Range("I6").Value = "Edoardo Zambanini"
$ws.Range("J6").Value = "Samuel Watson"
$ws.Range("K6").Value = "Anthon Charmig"
$ws.Range("L6").Value = "Andrea Raccagni Noviero"
$ws.Range("M6").Value = "Pierre Gautherat"

# Row 7 - Stage 5 results (2026-01-25)
$ws.Range("A7").Value = 46047
$ws.Range("B7").Value = "Santos Tour Down Under"
$ws.Range("C7").Value = "Stage 5"
$ws.Range("D7").Value = "Matthew Brennan"
$ws.Range("E7").Value = "Finn Fisher-Black"
$ws.Range("F7").Value = "Tobias Lund Andresen"
$ws.Range("G7").Value = "Brady Gilmore"
$ws.Range("H7").Value = "Simone Velasco"
$ws.Range("I7").Value = "Patrick Eddy"
$ws.Range("J7").Value = "Samuel Watson"
$ws.Range("K7").Value = "Edoardo Zambanini"
$ws.Range("L7").Value = "Natnael Tesfatsion"
$ws.Range("M7").Value = "Andrea Bagioli"

# Scroll the view so column E is the left-most visible column, and leave the
# final selection on H24:H25 (matches the saved workbook's view state).
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("H24:H25").Select()
